# ---------------------------------------------------------------------------
# Target diff analysis
# ---------------------------------------------------------------------------
# The supplied unified diff touches exactly four parts of the package:
#   word/document.xml, word/footer.xml, word/header.xml, word/styles.xml
#
# In every one of those files the only two kinds of change are:
#   1. The order in which the (same set of, ~68) xmlns:* namespace
#      declarations appear on the single root element
#      (<w:document>/<w:ftr>/<w:hdr>/<w:styles>) - no prefix/URI pair was
#      added, removed, or changed in value, they were merely re-emitted in
#      a different order.
#   2. In word/document.xml only, the text of a generator comment
#      ("<!-- Created by docx4j 11.4.9 ... on Linux -->" ->
#       "<!-- Created by docx4j 11.5.4 ... on Mac OS X -->").
#
# The accompanying commit message confirms this: it describes a source-level
# refactor of the authoring library itself (adding a `@DocxDsl` marker
# annotation and consolidating table-accessor helpers) and says nothing
# about changing this fixture's visible content. The XML diff we see is
# simply what happens when that library's test fixture gets regenerated by
# a newer docx4j release (11.4.9 -> 11.5.4) on a different build machine
# (Linux/Java 11 -> Mac OS X/Java 17): the serializer prints the same
# namespace set in a different (version-dependent) order and stamps a new
# generator comment. No paragraph text, run formatting, table, style
# definition, header/footer content, or any other reachable document
# object changed at all.
#
# Both of those two byte differences are artifacts of the XML *writer*
# that produced the file, not of the *document model* Word exposes via
# COM automation:
#   - Word (real or this interop shim) does not let a caller choose the
#     attribute-emission order of a part's root element; that is an
#     internal serialization detail, not a property on Document/Range/
#     Paragraph/etc.
#   - Freeform XML comments placed directly in a WordprocessingML part
#     (outside any w:t/w:instrText) are not part of the WordprocessingML
#     content model, so they are not visible to Find/Replace, Range.Text,
#     Range.XML, ContentControl.XML, etc. - there is no supported COM call
#     that can address/replace that one comment's characters, and in fact
#     the comment does not even survive a normal Word round trip (any
#     content edit forces a reserialization that silently drops it,
#     matching genuine Word behaviour for non-standard markup).
#
# Consequently there is no text/content/formatting edit to make here: the
# documents are identical at the object-model level before and after.
# Per the grading note "Only /tmp/work/edit.ps1 is graded", the safest and
# most faithful action is to leave the document exactly as-is (touch
# nothing), which also avoids incidentally introducing *new*, unrelated
# differences (e.g. this runtime adds a stray xmlns:w16du declaration and
# drops the generator comment as soon as *any* content mutation touches
# word/document.xml).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Touch the document object without mutating any content so the save is a
# faithful, byte-for-byte round trip of the (already up to date) body text,
# styles, header and footer.
$null = $d.Name
